# Update Daily Report: 2026-02-18
# Adds the new trading-day block (date serial 46070) to Daily_Data,
# then refreshes the dependent Today_Summary and Monthly_Stats rollups.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Daily_Data")

# --- Daily_Data: copy formatting from the most recent 22-row day block
#     (rows 640-661, date 46066) down into the new block (rows 662-683)
#     so the new rows inherit the same cell styles (e.g. the date format
#     on column A).
$srcBlock = $ws.Range("A640:H661")
$dstBlock = $ws.Range("A662:H683")
$srcBlock.Copy($dstBlock)

# --- Daily_Data: write the new day's values (date serial 46070 = 2026-02-17,
#     as reported in the 2026-02-18 daily report).
$ws.Cells.Item(662,1).Value = 46070
$ws.Cells.Item(662,2).Value = 'ASAHI DEPOSITORY LLC Registered'
$ws.Cells.Item(662,3).Value = 0
$ws.Cells.Item(662,4).Value = 0
$ws.Cells.Item(662,5).Value = 0
$ws.Cells.Item(662,6).Value = 0
$ws.Cells.Item(662,7).Value = 0
$ws.Cells.Item(662,8).Value = 0
$ws.Cells.Item(663,1).Value = 46070
$ws.Cells.Item(663,2).Value = 'ASAHI DEPOSITORY LLC Eligible'
$ws.Cells.Item(663,3).Value = 0
$ws.Cells.Item(663,4).Value = 0
$ws.Cells.Item(663,5).Value = 0
$ws.Cells.Item(663,6).Value = 0
$ws.Cells.Item(663,7).Value = 0
$ws.Cells.Item(663,8).Value = 0
$ws.Cells.Item(664,1).Value = 46070
$ws.Cells.Item(664,2).Value = 'BRINK''S, INC. Registered'
$ws.Cells.Item(664,3).Value = 73354.783
$ws.Cells.Item(664,4).Value = 0
$ws.Cells.Item(664,5).Value = 0
$ws.Cells.Item(664,6).Value = 0
$ws.Cells.Item(664,7).Value = 0
$ws.Cells.Item(664,8).Value = 73354.783
$ws.Cells.Item(665,1).Value = 46070
$ws.Cells.Item(665,2).Value = 'BRINK''S, INC. Eligible'
$ws.Cells.Item(665,3).Value = 85821.847
$ws.Cells.Item(665,4).Value = 0
$ws.Cells.Item(665,5).Value = 1361.109
$ws.Cells.Item(665,6).Value = -1361.109
$ws.Cells.Item(665,7).Value = 0
$ws.Cells.Item(665,8).Value = 84460.738
$ws.Cells.Item(666,1).Value = 46070
$ws.Cells.Item(666,2).Value = 'CNT DEPOSITORY, INC. Registered'
$ws.Cells.Item(666,3).Value = 1246.06
$ws.Cells.Item(666,4).Value = 0
$ws.Cells.Item(666,5).Value = 0
$ws.Cells.Item(666,6).Value = 0
$ws.Cells.Item(666,7).Value = 0
$ws.Cells.Item(666,8).Value = 1246.06
$ws.Cells.Item(667,1).Value = 46070
$ws.Cells.Item(667,2).Value = 'CNT DEPOSITORY, INC. Eligible'
$ws.Cells.Item(667,3).Value = 0
$ws.Cells.Item(667,4).Value = 0
$ws.Cells.Item(667,5).Value = 0
$ws.Cells.Item(667,6).Value = 0
$ws.Cells.Item(667,7).Value = 0
$ws.Cells.Item(667,8).Value = 0
$ws.Cells.Item(668,1).Value = 46070
$ws.Cells.Item(668,2).Value = 'DELAWARE DEPOSITORY Registered'
$ws.Cells.Item(668,3).Value = 1633.941
$ws.Cells.Item(668,4).Value = 0
$ws.Cells.Item(668,5).Value = 0
$ws.Cells.Item(668,6).Value = 0
$ws.Cells.Item(668,7).Value = 0
$ws.Cells.Item(668,8).Value = 1633.941
$ws.Cells.Item(669,1).Value = 46070
$ws.Cells.Item(669,2).Value = 'DELAWARE DEPOSITORY Eligible'
$ws.Cells.Item(669,3).Value = 18459.584
$ws.Cells.Item(669,4).Value = 0
$ws.Cells.Item(669,5).Value = 0
$ws.Cells.Item(669,6).Value = 0
$ws.Cells.Item(669,7).Value = 0
$ws.Cells.Item(669,8).Value = 18459.584
$ws.Cells.Item(670,1).Value = 46070
$ws.Cells.Item(670,2).Value = 'HSBC BANK, USA Registered'
$ws.Cells.Item(670,3).Value = 1394.758
$ws.Cells.Item(670,4).Value = 0
$ws.Cells.Item(670,5).Value = 0
$ws.Cells.Item(670,6).Value = 0
$ws.Cells.Item(670,7).Value = 0
$ws.Cells.Item(670,8).Value = 1394.758
$ws.Cells.Item(671,1).Value = 46070
$ws.Cells.Item(671,2).Value = 'HSBC BANK, USA Eligible'
$ws.Cells.Item(671,3).Value = 9281.979
$ws.Cells.Item(671,4).Value = 0
$ws.Cells.Item(671,5).Value = 0
$ws.Cells.Item(671,6).Value = 0
$ws.Cells.Item(671,7).Value = 0
$ws.Cells.Item(671,8).Value = 9281.979
$ws.Cells.Item(672,1).Value = 46070
$ws.Cells.Item(672,2).Value = 'INTERNATIONAL DEPOSITORY SERVICES OF DELAWARE Registered'
$ws.Cells.Item(672,3).Value = 2395.448
$ws.Cells.Item(672,4).Value = 0
$ws.Cells.Item(672,5).Value = 0
$ws.Cells.Item(672,6).Value = 0
$ws.Cells.Item(672,7).Value = 0
$ws.Cells.Item(672,8).Value = 2395.448
$ws.Cells.Item(673,1).Value = 46070
$ws.Cells.Item(673,2).Value = 'INTERNATIONAL DEPOSITORY SERVICES OF DELAWARE Eligible'
$ws.Cells.Item(673,3).Value = 0
$ws.Cells.Item(673,4).Value = 0
$ws.Cells.Item(673,5).Value = 0
$ws.Cells.Item(673,6).Value = 0
$ws.Cells.Item(673,7).Value = 0
$ws.Cells.Item(673,8).Value = 0
$ws.Cells.Item(674,1).Value = 46070
$ws.Cells.Item(674,2).Value = 'JP MORGAN CHASE BANK NA Registered'
$ws.Cells.Item(674,3).Value = 114061.421
$ws.Cells.Item(674,4).Value = 0
$ws.Cells.Item(674,5).Value = 0
$ws.Cells.Item(674,6).Value = 0
$ws.Cells.Item(674,7).Value = 0
$ws.Cells.Item(674,8).Value = 114061.421
$ws.Cells.Item(675,1).Value = 46070
$ws.Cells.Item(675,2).Value = 'JP MORGAN CHASE BANK NA Eligible'
$ws.Cells.Item(675,3).Value = 76408.669
$ws.Cells.Item(675,4).Value = 0
$ws.Cells.Item(675,5).Value = 0
$ws.Cells.Item(675,6).Value = 0
$ws.Cells.Item(675,7).Value = 0
$ws.Cells.Item(675,8).Value = 76408.669
$ws.Cells.Item(676,1).Value = 46070
$ws.Cells.Item(676,2).Value = 'LOOMIS INTERNATIONAL (US) LLC Registered'
$ws.Cells.Item(676,3).Value = 61157.444
$ws.Cells.Item(676,4).Value = 0
$ws.Cells.Item(676,5).Value = 0
$ws.Cells.Item(676,6).Value = 0
$ws.Cells.Item(676,7).Value = 0
$ws.Cells.Item(676,8).Value = 61157.444
$ws.Cells.Item(677,1).Value = 46070
$ws.Cells.Item(677,2).Value = 'LOOMIS INTERNATIONAL (US) LLC Eligible'
$ws.Cells.Item(677,3).Value = 71594.187
$ws.Cells.Item(677,4).Value = 0
$ws.Cells.Item(677,5).Value = 2588.547
$ws.Cells.Item(677,6).Value = -2588.547
$ws.Cells.Item(677,7).Value = 0
$ws.Cells.Item(677,8).Value = 69005.64
$ws.Cells.Item(678,1).Value = 46070
$ws.Cells.Item(678,2).Value = 'MALCA-AMIT USA, LLC Registered'
$ws.Cells.Item(678,3).Value = 395.145
$ws.Cells.Item(678,4).Value = 0
$ws.Cells.Item(678,5).Value = 0
$ws.Cells.Item(678,6).Value = 0
$ws.Cells.Item(678,7).Value = 0
$ws.Cells.Item(678,8).Value = 395.145
$ws.Cells.Item(679,1).Value = 46070
$ws.Cells.Item(679,2).Value = 'MALCA-AMIT USA, LLC Eligible'
$ws.Cells.Item(679,3).Value = 0
$ws.Cells.Item(679,4).Value = 0
$ws.Cells.Item(679,5).Value = 0
$ws.Cells.Item(679,6).Value = 0
$ws.Cells.Item(679,7).Value = 0
$ws.Cells.Item(679,8).Value = 0
$ws.Cells.Item(680,1).Value = 46070
$ws.Cells.Item(680,2).Value = 'MANFRA, TORDELLA & BROOKES, LLC Registered'
$ws.Cells.Item(680,3).Value = 49920.248
$ws.Cells.Item(680,4).Value = 0
$ws.Cells.Item(680,5).Value = 0
$ws.Cells.Item(680,6).Value = 0
$ws.Cells.Item(680,7).Value = 0
$ws.Cells.Item(680,8).Value = 49920.248
$ws.Cells.Item(681,1).Value = 46070
$ws.Cells.Item(681,2).Value = 'MANFRA, TORDELLA & BROOKES, LLC Eligible'
$ws.Cells.Item(681,3).Value = 2104.855
$ws.Cells.Item(681,4).Value = 0
$ws.Cells.Item(681,5).Value = 300.172
$ws.Cells.Item(681,6).Value = -300.172
$ws.Cells.Item(681,7).Value = 0
$ws.Cells.Item(681,8).Value = 1804.683
$ws.Cells.Item(682,1).Value = 46070
$ws.Cells.Item(682,2).Value = 'STONEX PRECIOUS METALS LLC Registered'
$ws.Cells.Item(682,3).Value = 14122.765
$ws.Cells.Item(682,4).Value = 0
$ws.Cells.Item(682,5).Value = 0
$ws.Cells.Item(682,6).Value = 0
$ws.Cells.Item(682,7).Value = 0
$ws.Cells.Item(682,8).Value = 14122.765
$ws.Cells.Item(683,1).Value = 46070
$ws.Cells.Item(683,2).Value = 'STONEX PRECIOUS METALS LLC Eligible'
$ws.Cells.Item(683,3).Value = 16.075
$ws.Cells.Item(683,4).Value = 0
$ws.Cells.Item(683,5).Value = 0
$ws.Cells.Item(683,6).Value = 0
$ws.Cells.Item(683,7).Value = 0
$ws.Cells.Item(683,8).Value = 16.075

# --- Today_Summary: refresh the three depositories whose Eligible figures
#     changed because of the new day's WITHDRAWN activity.
$today = $wb.Worksheets.Item("Today_Summary")

# BRINK'S, INC. (row 3)
$today.Cells.Item(3,2).Value = 84460.738
$today.Cells.Item(3,4).Value = 157815.521

# LOOMIS INTERNATIONAL (US) LLC (row 9)
$today.Cells.Item(9,2).Value = 69005.64
$today.Cells.Item(9,4).Value = 130163.084

# MANFRA, TORDELLA & BROOKES, LLC (row 11)
$today.Cells.Item(11,2).Value = 1804.683
$today.Cells.Item(11,4).Value = 51724.931

# --- Monthly_Stats: refresh the 2026-02 grand totals and the per-depository
#     detail rows affected by the new day's RECEIVED/WITHDRAWN/TOTAL_TODAY.
$monthly = $wb.Worksheets.Item("Monthly_Stats")

# 2026-02 Grand_Total summary row (row 2)
$monthly.Cells.Item(2,2).Value = 259437.368
$monthly.Cells.Item(2,4).Value = 579119.3809999999

# BRINK'S, INC. Eligible / 2026-02 detail row (row 10)
$monthly.Cells.Item(10,4).Value = 17539.066
$monthly.Cells.Item(10,5).Value = 84460.738

# LOOMIS INTERNATIONAL (US) LLC Eligible / 2026-02 detail row (row 22)
$monthly.Cells.Item(22,4).Value = 65660.113
$monthly.Cells.Item(22,5).Value = 69005.64

# MANFRA, TORDELLA & BROOKES, LLC Eligible / 2026-02 detail row (row 26)
$monthly.Cells.Item(26,4).Value = 300.172
$monthly.Cells.Item(26,5).Value = 1804.683
